# "Schedule/每日做题计划.xlsx" — add a new row for problem 304 ("dxs: problem 304")
#
# Target state (row 27 on the active "新题" sheet):
#   A27 = 3/30/2019  (one day after A26's 3/29/2019, same date style)
#   B27 = 304
#   F27 = "done"
# Dimension grows to A1:F27 and the selection moves to F27.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy A26 into A27 first so A27 inherits the exact same date number-format
# style (xf index) used by the rest of column A, rather than Excel inventing
# a brand-new number-format style for a bare date value/string.
$ws.Range("A26").Copy($ws.Range("A27"))
# Now overwrite with the actual date for this row: 2019-03-30 == serial 43554.
$ws.Range("A27").Value = 43554

# New problem number solved.
$ws.Range("B27").Value = 304

# Status column.
$ws.Range("F27").Value = "done"

# Match the workbook's saved cursor position after the edit.
$ws.Range("F27").Select()
